$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets("Lookup")

# --- Add "Rapid Evac" to the Evac column (K) for the existing HACE / HAPE rows ---
$ws.Cells.Item(20, 11).Value = "Rapid Evac"
$ws.Cells.Item(21, 11).Value = "Rapid Evac"

# --- New condition rows appended below the existing table (rows 22-46) ---
# Each entry is: row, Name(A), Symptom Keywords(I), Treatments(J), Evac(K)
$rows = @(
  @{ r = 22; a = "Myocardial Infarcation"; i = "Lightheadedness, dizziness, anxiety, denial, radiating pain, persistent chest pain, "; j = "Reduce anxiety and activity, Place patient in position of comfort, support patient with personal medications"; k = "Evac all patients with suspected cardiac chset pain, rapid evac any patient with a new onset of chest pain that is not clearly musculoskeletal" },
  @{ r = 23; a = "Hyperventilation" },
  @{ r = 24; a = "Pulmonary Emboism" },
  @{ r = 25; a = "Pnuemonia"; i = "productive cough, wet lung sounds, malaise, fatigue, shortness of breath" },
  @{ r = 26; a = "Asthma" },
  @{ r = 27; a = "Stroke" },
  @{ r = 28; a = "Seizure" },
  @{ r = 29; a = "Syncope" },
  @{ r = 30; a = "Gastroenteritis" },
  @{ r = 31; a = "Mild/Moderate Allergic Reactions" },
  @{ r = 32; a = "Anaphylaxis" },
  @{ r = 33; a = "Hyperglycemia" },
  @{ r = 34; a = "Hypoglycemia" },
  @{ r = 35; a = "Painful Menstruation" },
  @{ r = 36; a = "Ectopic Pregnancy" },
  @{ r = 37; a = "Vaginitis" },
  @{ r = 38; a = "UTI" },
  @{ r = 39; a = "Testicular Torsion" },
  @{ r = 40; a = "Epididymitis" },
  @{ r = 41; a = "Inguinal Hernia" },
  @{ r = 42; a = "Stress Injury" },
  @{ r = 43; a = "Anxiety" },
  @{ r = 44; a = "Depression" },
  @{ r = 45; a = "Mania/Psychosis" },
  @{ r = 46; a = "Suicidal Thoughts" }
)

foreach ($row in $rows) {
  $ws.Cells.Item($row.r, 1).Value = $row.a
  if ($row.ContainsKey("i")) { $ws.Cells.Item($row.r, 9).Value = $row.i }
  if ($row.ContainsKey("j")) { $ws.Cells.Item($row.r, 10).Value = $row.j }
  if ($row.ContainsKey("k")) { $ws.Cells.Item($row.r, 11).Value = $row.k }
}

# --- Match the author's final on-screen state: scrolled down, new cell selected ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 15
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A47").Select() | Out-Null
